$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 346 ("「一人ひとりが責任を」" Covid post) which removes the row
# and shifts all subsequent rows up by one.
$ws.Rows.Item(346).Delete()
